$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 1 data values (columns C..N) ---
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 8
$ws.Range("F1").Value = 31
$ws.Range("G1").Value = 12
$ws.Range("H1").Value = 19
$ws.Range("I1").Value = 24
$ws.Range("J1").Value = 23
$ws.Range("K1").Value = 0.097000000000000003
$ws.Range("L1").Value = 0.024
$ws.Range("M1").Value = 0.054000000000000006
$ws.Range("N1").Value = 0.094

# --- Update column widths to match the new content widths ---
# (columns shrink from 2-digit to 1-digit numbers, or grow from 1-digit to
# 2-digit numbers, so their widths are swapped accordingly)
$ws.Columns.Item(3).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(4).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(5).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(9).ColumnWidth = 2.3333333333333335
